$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 10; this shifts the existing
# rows 10-14 down to 11-15 (matching the diff's row-by-row shift).
$ws.Rows("10").Insert()

# Populate the newly inserted row 10 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,J,K,L,R mirror the record that used to sit in
# row 10 (now row 11); only the date/volume/price/unit/kg columns differ.
$ws.Cells.Item(10, 1).Value = 1
$ws.Cells.Item(10, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(10, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(10, 4).Value = "10/05/2022"
$ws.Cells.Item(10, 5).Value = 15
$ws.Cells.Item(10, 6).Value = "Fruta"
$ws.Cells.Item(10, 7).Value = 100107
$ws.Cells.Item(10, 8).Value = "Otros"
$ws.Cells.Item(10, 9).Value = 100107002
$ws.Cells.Item(10, 10).Value = "Chirimoya"
$ws.Cells.Item(10, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(10, 12).Value = "Segunda"
$ws.Cells.Item(10, 13).Value = 160
$ws.Cells.Item(10, 14).Value = 26000
$ws.Cells.Item(10, 15).Value = 27000
$ws.Cells.Item(10, 16).Value = 26500
$ws.Cells.Item(10, 17).Value = "$/caja 12 kilos"
$ws.Cells.Item(10, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(10, 19).Value = 2208
$ws.Cells.Item(10, 20).Value = 12
